$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row, column index, new text value.
# Values are written via a temporary "@" (Text) NumberFormat so that
# numeric-looking strings (e.g. "1.00") are stored as text, matching the
# original inline-string cell type; the cell's original Style is restored
# afterward so no visible formatting changes.
$updates = @(
    ,@(2, 4, '66.718.14')
    ,@(2, 5, '  +1.98%  ')
    ,@(3, 4, '3.696.86')
    ,@(3, 5, '  +4.56%  ')
    ,@(4, 4, '1.00')
    ,@(4, 5, '  +0.29%  ')
    ,@(5, 4, '419.23')
    ,@(5, 5, '  -0.96%  ')
    ,@(6, 4, '130.12')
    ,@(6, 5, '  -1.31%  ')
    ,@(7, 4, '3.690.91')
    ,@(7, 5, '  +4.58%  ')
    ,@(8, 4, '0.641')
    ,@(8, 5, '  -0.13%  ')
    ,@(9, 5, '  -0.10%  ')
    ,@(10, 5, '  -3.65%  ')
    ,@(11, 4, '0.180')
    ,@(11, 5, '  +6.05%  ')
    ,@(12, 4, '0.0000391')
    ,@(12, 5, '  +42.20%  ')
    ,@(13, 4, '43.15')
    ,@(13, 5, '  -1.12%  ')
    ,@(14, 4, '10.65')
    ,@(14, 5, '  +5.25%  ')
    ,@(15, 4, '4.279.56')
    ,@(15, 5, '  +4.90%  ')
    ,@(16, 4, '0.140')
    ,@(16, 5, '  -0.88%  ')
    ,@(17, 4, '3.810.62')
    ,@(17, 5, '  +7.99%  ')
    ,@(18, 4, '20.49')
    ,@(18, 5, '  -1.30%  ')
    ,@(19, 5, '  +5.01%  ')
    ,@(20, 5, '  +0.80%  ')
    ,@(21, 4, '66.745.67')
    ,@(21, 5, '  +2.38%  ')
    ,@(22, 4, '440.41')
    ,@(22, 5, '  -5.90%  ')
    ,@(23, 4, '16.43')
    ,@(23, 5, '  +20.61%  ')
    ,@(24, 4, '89.56')
    ,@(24, 5, '  -2.57%  ')
    ,@(25, 4, '3.14')
    ,@(25, 5, '  -4.83%  ')
    ,@(26, 4, '37.46')
    ,@(26, 5, '  +8.23%  ')
    ,@(27, 4, '10.32')
    ,@(27, 5, '  +0.80%  ')
    ,@(28, 5, '  -2.04%  ')
    ,@(29, 5, '  +3.94%  ')
    ,@(30, 4, '0.125')
    ,@(30, 5, '  +8.97%  ')
    ,@(31, 4, '12.79')
    ,@(31, 5, '  +1.82%  ')
    ,@(32, 4, '2.78')
    ,@(32, 5, '  +2.45%  ')
    ,@(33, 4, '7.27')
    ,@(33, 5, '  -4.61%  ')
    ,@(34, 4, '0.166')
    ,@(34, 5, '  -1.20%  ')
    ,@(35, 4, '41.76')
    ,@(35, 5, '  +2.35%  ')
    ,@(36, 4, '57.07')
    ,@(36, 5, '  -1.56%  ')
    ,@(37, 5, '  +0.04%  ')
    ,@(38, 5, '  -4.82%  ')
    ,@(39, 4, '3.09')
    ,@(39, 5, '  +32.54%  ')
    ,@(40, 4, '0.0₃0732')
    ,@(40, 5, '  +2.10%  ')
    ,@(41, 4, '0.150')
    ,@(41, 5, '  +4.41%  ')
    ,@(42, 4, '28.51')
    ,@(42, 5, '  +28.69%  ')
    ,@(43, 4, '0.998')
    ,@(43, 5, '  +0.19%  ')
    ,@(44, 5, '  +0.48%  ')
    ,@(45, 4, '148.93')
    ,@(45, 5, '  +1.75%  ')
    ,@(46, 5, '  +3.54%  ')
    ,@(47, 2, 'Stacks')
    ,@(47, 3, 'https://coinranking.com/coin/mMPrMcB7+stacks-stx')
    ,@(47, 4, '2.88')
    ,@(47, 5, '  -8.85%  ')
    ,@(48, 2, 'NEARProtocol')
    ,@(48, 3, 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near')
    ,@(48, 4, '4.35')
    ,@(48, 5, '  -4.31%  ')
    ,@(49, 4, '2.60')
    ,@(49, 5, '  -7.30%  ')
    ,@(50, 5, '  -5.52%  ')
    ,@(51, 5, '  +11.23%  ')
)

foreach ($u in $updates) {
    $r = $u[0]; $c = $u[1]; $val = $u[2]
    $cell = $ws.Cells.Item($r, $c)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}
